# ---------------------------------------------------------------------------
# Ebates-COLING-2016-results.xlsx edit
#   * add a new worksheet "Earlier" (after KLvsPrediction) holding the
#     "earlier" LogReg / XGboost micro precision+F1 numbers, with a bar
#     chart comparing them
#   * drop a "92 is a placeholder" note into Predictions!G27
#   * update the saved selections
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Predictions!G27 placeholder note
# ---------------------------------------------------------------------------
$predictions = $wb.Worksheets.Item("Predictions")
$predictions.Range("G27").Value = "92 is a placeholder"

# ---------------------------------------------------------------------------
# 2. New "Earlier" worksheet, inserted after the last existing sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$earlier = $wb.Worksheets.Add($null, $lastSheet)
$earlier.Name = "Earlier"

# Header row
$headers = @("Category", "LogReg L1 Micro Precision", "LogReg L1 Micro F1", "XGboost Micro Precision", "XGboost Micro F1")
for ($c = 1; $c -le 5; $c++) {
    $cell = $earlier.Cells.Item(2, $c)
    $cell.Value = $headers[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

# Data rows (category name, LogReg precision, LogReg F1, XGboost precision, XGboost F1)
$rows = @(
    @("Toys", 76.003999999999991, 76.003999999999991, 81.347999999999999, 81.347999999999999),
    @("Home, Patio and Furniture", 84.00200000000001, 84.00200000000001, 88.9, 88.8),
    @("Jewelry and Watches", 80.894000000000005, 80.894000000000005, 86.027999999999992, 86.027999999999992),
    @("Bags, Handbags and Accessories", 81.676000000000002, 81.676000000000002, 85.001999999999995, 85.001999999999995),
    @("Health, Beauty and Fragrance", 82.177999999999997, 82.177999999999997, 85.390000000000015, 85.390000000000015),
    @("Shoes", 64.516000000000005, 64.516000000000005, 66.608000000000004, 66.608000000000004),
    @("Electronics and Computers", 80.132000000000005, 80.132000000000005, 84.496000000000009, 84.496000000000009),
    @("Office", 89.217999999999989, 89.217999999999989, 92.891999999999996, 92.891999999999996),
    @("Sports and Fitness", 83.571999999999989, 83.571999999999989, 87.443999999999988, 87.443999999999988),
    @("Automotive", 88.794000000000011, 88.794000000000011, 94.753999999999991, 94.753999999999991),
    @("Industrial", 88.236000000000018, 88.236000000000018, 93.137999999999991, 93.137999999999991),
    @("Baby Products", 88.190000000000012, 88.190000000000012, 89.671999999999997, 89.671999999999997),
    @("Baby and Kids Clothes", 89.033999999999992, 89.033999999999992, 92.140000000000015, 92.140000000000015),
    @("Men's Clothing", 82.419999999999987, 82.419999999999987, 85.464000000000013, 85.464000000000013),
    @("Women's Clothing", 83.361999999999995, 83.361999999999995, 85.468000000000004, 85.468000000000004)
)

$r = 3
foreach ($row in $rows) {
    $earlier.Cells.Item($r, 1).Value = $row[0]
    $earlier.Cells.Item($r, 2).Value = $row[1]
    $earlier.Cells.Item($r, 3).Value = $row[2]
    $earlier.Cells.Item($r, 4).Value = $row[3]
    $earlier.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

$earlier.Columns.Item(1).ColumnWidth = 27.33203125

# ---------------------------------------------------------------------------
# 3. Bar chart on the "Earlier" sheet (precision + F1 for the two "Level 1"
#    classifiers, LogReg L1 vs XGboost)
# ---------------------------------------------------------------------------
$chartObj = $earlier.ChartObjects().Add(150, 20, 520, 330)
$chart = $chartObj.Chart
$chart.ChartType = 51
$chart.SetSourceData($earlier.Range("A2:C17"))

$chart.HasTitle = $true
$chart.ChartTitle.Text = "State-of-the-art Level 1 Classifier MICRO Performance Comparison"

$chart.HasLegend = $true
$chart.Legend.Position = -4152

$valueAxis = $chart.Axes(2)
$valueAxis.MinimumScale = 60
$valueAxis.MaximumScale = 100
$valueAxis.HasTitle = $true
$valueAxis.AxisTitle.Text = "Prediction Micro Precision/F1"

# ---------------------------------------------------------------------------
# 4. Selections / active sheet bookkeeping
# ---------------------------------------------------------------------------
$earlier.Range("B3").Select()
$predictions.Activate()
$predictions.Range("G28").Select()

Write-Host "edit applied"
